# Scheduled runner update: refresh cached market-price / profit figures
# (columns H-N: currentAveragePrice*, LevePrice*, LeveProfit*) across the
# ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1180
$ws.Range("I32").Value = 800
$ws.Range("J32").Value = 1750
$ws.Range("K32").Value = 800
$ws.Range("L32").Value = 1750
$ws.Range("M32").Value = -474
$ws.Range("N32").Value = -2402
$ws.Range("H53").Value = 3172
$ws.Range("I53").Value = 489
$ws.Range("J53").Value = 10326.667
$ws.Range("K53").Value = 489
$ws.Range("L53").Value = 10326.667
$ws.Range("M53").Value = 148
$ws.Range("N53").Value = -11600.667
$ws.Range("H62").Value = 3523.1155
$ws.Range("I62").Value = 2770
$ws.Range("K62").Value = 2770
$ws.Range("M62").Value = -2146
$ws.Range("H65").Value = 3523.1155
$ws.Range("I65").Value = 2770
$ws.Range("K65").Value = 13850
$ws.Range("M65").Value = -10730
$ws.Range("H76").Value = 4632454.5
$ws.Range("I76").Value = 3150
$ws.Range("J76").Value = 6947107
$ws.Range("K76").Value = 3150
$ws.Range("L76").Value = 6947107
$ws.Range("M76").Value = -2835
$ws.Range("N76").Value = -6947737
$ws.Range("H79").Value = 4632454.5
$ws.Range("I79").Value = 3150
$ws.Range("J79").Value = 6947107
$ws.Range("K79").Value = 3150
$ws.Range("L79").Value = 6947107
$ws.Range("M79").Value = -2058
$ws.Range("N79").Value = -6949291
$ws.Range("H86").Value = 9816
$ws.Range("I86").Value = 1550
$ws.Range("K86").Value = 1550
$ws.Range("M86").Value = -427
$ws.Range("H89").Value = 9816
$ws.Range("I89").Value = 1550
$ws.Range("K89").Value = 7750
$ws.Range("M89").Value = -2134
$ws.Range("H112").Value = 1126.8276
$ws.Range("J112").Value = 1151.1111
$ws.Range("L112").Value = 3453.3333
$ws.Range("N112").Value = -5669.3333
$ws.Range("H113").Value = 83337576
$ws.Range("I113").Value = 200001520
$ws.Range("J113").Value = 6182.857
$ws.Range("K113").Value = 200001520
$ws.Range("L113").Value = 6182.857
$ws.Range("M113").Value = -199998266
$ws.Range("N113").Value = -12690.857
$ws.Range("H125").Value = 349.2857
$ws.Range("I125").Value = 309.14285
$ws.Range("K125").Value = 2782.28565
$ws.Range("M125").Value = -322.2856500000003
$ws.Range("H132").Value = 2101.0222
$ws.Range("I132").Value = 2101.0222
$ws.Range("K132").Value = 6303.0666
$ws.Range("M132").Value = -3773.0666
$ws.Range("H137").Value = 65967.87
$ws.Range("I137").Value = 82606.19
$ws.Range("J137").Value = 3254.2307
$ws.Range("K137").Value = 247818.57
$ws.Range("L137").Value = 9762.6921
$ws.Range("M137").Value = -245268.57
$ws.Range("N137").Value = -14862.6921

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1775.5217
$ws.Range("I61").Value = 1433.5366
$ws.Range("J61").Value = 4579.8
$ws.Range("K61").Value = 1433.5366
$ws.Range("L61").Value = 4579.8
$ws.Range("M61").Value = -1221.5366
$ws.Range("N61").Value = -5003.8
$ws.Range("H63").Value = 2843351.8
$ws.Range("I63").Value = 2687
$ws.Range("J63").Value = 31250000
$ws.Range("K63").Value = 2687
$ws.Range("L63").Value = 31250000
$ws.Range("M63").Value = -2001
$ws.Range("N63").Value = -31251372
$ws.Range("H66").Value = 2843351.8
$ws.Range("I66").Value = 2687
$ws.Range("J66").Value = 31250000
$ws.Range("K66").Value = 13435
$ws.Range("L66").Value = 156250000
$ws.Range("M66").Value = -10003
$ws.Range("N66").Value = -156256864
$ws.Range("H132").Value = 9894.403
$ws.Range("I132").Value = 1646.7954
$ws.Range("J132").Value = 30055.223
$ws.Range("K132").Value = 4940.3862
$ws.Range("L132").Value = 90165.66900000001
$ws.Range("M132").Value = -2410.3862
$ws.Range("N132").Value = -95225.66900000001
$ws.Range("H136").Value = 1775.5217
$ws.Range("I136").Value = 1433.5366
$ws.Range("J136").Value = 4579.8
$ws.Range("K136").Value = 4300.6098
$ws.Range("L136").Value = 13739.4
$ws.Range("M136").Value = -1750.6098
$ws.Range("N136").Value = -18839.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1465.5385
$ws.Range("I99").Value = 1392.8572
$ws.Range("K99").Value = 1392.8572
$ws.Range("M99").Value = 105.1428000000001
$ws.Range("H105").Value = 1472676.1
$ws.Range("I105").Value = 1444.4445
$ws.Range("J105").Value = 2002319.5
$ws.Range("K105").Value = 1444.4445
$ws.Range("L105").Value = 2002319.5
$ws.Range("M105").Value = 302.5554999999999
$ws.Range("N105").Value = -2005813.5
$ws.Range("H134").Value = 3680
$ws.Range("I134").Value = 3664.1025
$ws.Range("J134").Value = 3990
$ws.Range("K134").Value = 10992.3075
$ws.Range("L134").Value = 11970
$ws.Range("M134").Value = -8457.307499999999
$ws.Range("N134").Value = -17040

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 881.9375
$ws.Range("I16").Value = 788.36365
$ws.Range("J16").Value = 1087.8
$ws.Range("K16").Value = 788.36365
$ws.Range("L16").Value = 1087.8
$ws.Range("M16").Value = -501.36365
$ws.Range("N16").Value = -1661.8
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("H113").Value = 881.9375
$ws.Range("I113").Value = 788.36365
$ws.Range("J113").Value = 1087.8
$ws.Range("K113").Value = 788.36365
$ws.Range("L113").Value = 1087.8
$ws.Range("M113").Value = 1381.63635
$ws.Range("N113").Value = -5427.8
$ws.Range("H134").Value = 1246.2632
$ws.Range("I134").Value = 941.5
$ws.Range("K134").Value = 2824.5
$ws.Range("M134").Value = -289.5
$ws.Range("M38").ClearContents()
$ws.Range("M46").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 137.42857
$ws.Range("J12").Value = 234
$ws.Range("L12").Value = 702
$ws.Range("N12").Value = -1048
$ws.Range("H38").Value = 81.8
$ws.Range("I38").Value = 81.8
$ws.Range("K38").Value = 245.4
$ws.Range("M38").Value = 101.6
$ws.Range("H107").Value = 7951.6
$ws.Range("I107").Value = 10839
$ws.Range("K107").Value = 32517
$ws.Range("M107").Value = -30597
$ws.Range("H131").Value = 634.1799999999999
$ws.Range("I131").Value = 307.7857
$ws.Range("J131").Value = 761.1111
$ws.Range("K131").Value = 923.3571000000001
$ws.Range("L131").Value = 2283.3333
$ws.Range("M131").Value = 4116.6429
$ws.Range("N131").Value = -12363.3333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4177544
$ws.Range("I70").Value = 5350
$ws.Range("K70").Value = 5350
$ws.Range("M70").Value = -5080
$ws.Range("H73").Value = 4177544
$ws.Range("I73").Value = 5350
$ws.Range("K73").Value = 5350
$ws.Range("M73").Value = -4414
$ws.Range("H80").Value = 17860554
$ws.Range("I80").Value = 31253020
$ws.Range("J80").Value = 3933.3333
$ws.Range("K80").Value = 31253020
$ws.Range("L80").Value = 3933.3333
$ws.Range("M80").Value = -31252022
$ws.Range("N80").Value = -5929.3333
$ws.Range("H83").Value = 17860554
$ws.Range("I83").Value = 31253020
$ws.Range("J83").Value = 3933.3333
$ws.Range("K83").Value = 156265100
$ws.Range("L83").Value = 19666.6665
$ws.Range("M83").Value = -156260108
$ws.Range("N83").Value = -29650.6665
$ws.Range("H97").Value = 2644.4211
$ws.Range("I97").Value = 2471.5
$ws.Range("K97").Value = 2471.5
$ws.Range("M97").Value = -1975.5
$ws.Range("H107").Value = 290.5625
$ws.Range("I107").Value = 289
$ws.Range("J107").Value = 293.16666
$ws.Range("K107").Value = 289
$ws.Range("L107").Value = 293.16666
$ws.Range("M107").Value = 1631
$ws.Range("N107").Value = -4133.16666
$ws.Range("H113").Value = 10301.1
$ws.Range("I113").Value = 13673
$ws.Range("J113").Value = 2433.3333
$ws.Range("K113").Value = 13673
$ws.Range("L113").Value = 2433.3333
$ws.Range("M113").Value = -11503
$ws.Range("N113").Value = -6773.3333
$ws.Range("H132").Value = 12136.272
$ws.Range("I132").Value = 3244.7812
$ws.Range("K132").Value = 9734.3436
$ws.Range("M132").Value = -7204.3436

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 612.5909
$ws.Range("I16").Value = 646
$ws.Range("J16").Value = 564.3333
$ws.Range("K16").Value = 646
$ws.Range("L16").Value = 564.3333
$ws.Range("M16").Value = -476
$ws.Range("N16").Value = -904.3333
$ws.Range("H22").Value = 3625.25
$ws.Range("J22").Value = 1500
$ws.Range("L22").Value = 1500
$ws.Range("N22").Value = -2090
$ws.Range("H27").Value = 3625.25
$ws.Range("J27").Value = 1500
$ws.Range("L27").Value = 1500
$ws.Range("N27").Value = -1714
$ws.Range("H40").Value = 3962.2727
$ws.Range("I40").Value = 3609.25
$ws.Range("K40").Value = 3609.25
$ws.Range("M40").Value = -3473.25
$ws.Range("H93").Value = 2514.389
$ws.Range("I93").Value = 2417.2666
$ws.Range("J93").Value = 3000
$ws.Range("K93").Value = 2417.2666
$ws.Range("L93").Value = 3000
$ws.Range("M93").Value = -1169.2666
$ws.Range("N93").Value = -5496
$ws.Range("H132").Value = 229565.7
$ws.Range("I132").Value = 318783.62
$ws.Range("J132").Value = 3546.9333
$ws.Range("K132").Value = 956350.86
$ws.Range("L132").Value = 10640.7999
$ws.Range("M132").Value = -953820.86
$ws.Range("N132").Value = -15700.7999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 797.4737
$ws.Range("I132").Value = 567.7646999999999
$ws.Range("J132").Value = 2750
$ws.Range("K132").Value = 1703.2941
$ws.Range("L132").Value = 8250
$ws.Range("M132").Value = 826.7059000000002
$ws.Range("N132").Value = -13310
$ws.Range("H136").Value = 28676708
$ws.Range("I136").Value = 39703640
$ws.Range("J136").Value = 6679.8
$ws.Range("K136").Value = 119110920
$ws.Range("L136").Value = 20039.4
$ws.Range("M136").Value = -119108370
$ws.Range("N136").Value = -25139.4
